$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1/Q1 (continue the 0..15 sequence from row 1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header formatting (bold font, border, centered alignment) from O1
# onto the two new header cells without disturbing their values.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rows 2-25: swap the I/K and M/O column values, and append two new
# columns P and Q (both value 2) for every data row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
